# Updated Last Week of Standups
# Milestone3.2 sheet: Dylan (row 2) logged 1 hour on J2 (8-Nov) and 3 hours
# on M2 (11-Nov). Dependent formulas (U2, J10, M10, J11:T11) and the
# burndown chart's cached series recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Milestone3.2")
$ws.Activate()

$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 3

# Match the author's final selection on that sheet.
$ws.Range("O25").Select()
